$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap "Ano" and "Variável" between columns B and C ---
$ws.Range("B1").Value = "Ano"
$ws.Range("C1").Value = "Variável"

# --- Data rows: A=Região, B=Ano (date text), C=Variável, D=Valor, E=Posição relativamente às demais UF ---
$variavel = "Roubo seguido de morte (latrocínio)"

# Região, Ano, Valor, Posição (blank for Brasil/Nordeste rows)
$data = @(
    @("Brasil",   "01/01/2017", 1.490280650842423,  $null),
    @("Brasil",   "01/01/2018", 1.292246920939238,  $null),
    @("Brasil",   "01/01/2019", 1.050170774581274,  $null),
    @("Brasil",   "01/01/2020", 0.9528706417276306, $null),
    @("Brasil",   "01/01/2021", 0.9417499517438702, $null),
    @("Brasil",   "01/01/2022", 0.7697804286001625, $null),
    @("Brasil",   "01/01/2023", 0.5941299031652268, $null),
    @("Brasil",   "01/01/2024", 0.455588497772493,  $null),
    @("Nordeste", "01/01/2017", 1.670186678902145,  $null),
    @("Nordeste", "01/01/2018", 1.291955013741005,  $null),
    @("Nordeste", "01/01/2019", 1.009950045575144,  $null),
    @("Nordeste", "01/01/2020", 1.028997954630678,  $null),
    @("Nordeste", "01/01/2021", 0.9736065861473615, $null),
    @("Nordeste", "01/01/2022", 0.7913429844997724, $null),
    @("Nordeste", "01/01/2023", 0.5703241188109771, $null),
    @("Nordeste", "01/01/2024", 0.4934573726027487, $null),
    @("Sergipe",  "01/01/2017", 2.569480070137946,  5),
    @("Sergipe",  "01/01/2018", 1.404551096691053,  9),
    @("Sergipe",  "01/01/2019", 1.087475673169191,  11),
    @("Sergipe",  "01/01/2020", 1.250521769427934,  8),
    @("Sergipe",  "01/01/2021", 0.7269032893656614, 16),
    @("Sergipe",  "01/01/2022", 0.4665286872730974, 19),
    @("Sergipe",  "01/01/2023", 0.25247775355394,   24),
    @("Sergipe",  "01/01/2024", 0.2923385581194111, 20)
)

# Column B holds dd/mm/yyyy-looking text; force text format so Excel
# doesn't silently convert the literal strings into date serials.
$ws.Range("B2:B25").NumberFormat = "@"

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $variavel
    $ws.Cells.Item($row, 4).Value = $rec[2]

    $posicao = $rec[3]
    if ($posicao -ne $null) {
        $ws.Cells.Item($row, 5).Value = $posicao
    } else {
        $ws.Cells.Item($row, 5).ClearContents()
    }
    $row++
}

Write-Output "Applied g19.3 data refresh"
